$d = $word.ActiveDocument

# Remove the entire bullet paragraph "There will be a selector to add new
# words, and one to delete a word" (including its own paragraph mark),
# leaving the preceding "bubble sort" paragraph and the following
# paragraph ("The program will detect...") fully intact with their own
# original paragraph marks/properties.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "There will be a selector to add new words*") {
        $p.Range.Delete()
        break
    }
}
